$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume change (E) columns for rows 2-51
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.510.69"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.456.91"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.69"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.03"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  +4.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.41"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.838.21"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.84"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.80"
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.431.65"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.770"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.540.77"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.73"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.32"
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.21"
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.70"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.24"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.03"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.09"
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0759"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.38"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.114"
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.92"
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.971.71"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0281"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.59"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.93"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.696.63"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.30"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.34"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.171"
$ws.Range("E51").Value = "  -1.20%  "
